# Initial Data File Updated
# Adds two new transaction rows (92-93, "Pago de Axtel" / "Axtel") to the
# "Transacciones" sheet, right after the existing last row (91), keeping
# the running-balance formulas (L, N, O, P) consistent with the rest of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# Insert two fresh rows right after the current last data row (91); Excel
# carries the formatting of row 91 down onto the new rows 92:93.
$ws.Range("A92:A93").EntireRow.Insert()

# --- Row 92: "Pago de Axtel" (Servicios / Gasto, Tarjeta Santander) ---
$ws.Range("A92").Value2 = 43578
$ws.Range("B92").Value2 = 429
$ws.Range("C92").Value2 = "Pago de Axtel"
$ws.Range("D92").Value2 = "Servicios"
$ws.Range("E92").Value2 = "Gasto"
$ws.Range("F92").Value2 = "Tarjeta Santander"
$ws.Range("G92").Value2 = "Axtel"
$ws.Range("K92").Value2 = 6769.44
$ws.Range("L92").Formula = "=L91-B92"
$ws.Range("M92").Value2 = 2

# --- Row 93: "Pago de Axtel" (Pagos / Ingreso, Tarjeta Santander) ---
$ws.Range("A93").Value2 = 43578
$ws.Range("B93").Value2 = 285
$ws.Range("C93").Value2 = "Pago de Axtel"
$ws.Range("D93").Value2 = "Pagos"
$ws.Range("E93").Value2 = "Ingreso"
$ws.Range("F93").Value2 = "Tarjeta Santander"
$ws.Range("G93").Value2 = "NA"
$ws.Range("K93").Value2 = 6769.44
$ws.Range("L93").Formula = "=L92+B93"
$ws.Range("M93").Value2 = 2

# Running-total formulas, filled together so the pair shares one formula
# group the same way the rest of the sheet does.
$ws.Range("N92:N93").Formula = "=SUM(K92:M92)"
$ws.Range("O92:O93").Formula = "=N92-4000"
$ws.Range("P92").Formula = "=O92-Ahorros!`$E`$4"
$ws.Range("P93").Formula = "=O93-Ahorros!`$E`$4"

# Keep the view close to where it was (bottom of the table, last new cell
# selected) like the source workbook shows after the edit.
$ws.Activate()
$ws.Range("P93").Select()

$wb.Application.Calculate()
